# Edit script: applies the resume content revision described by the diff.
$d  = $word.ActiveDocument
$cr = [string][char]13
$bu = [string][char]8226

# ---------------------------------------------------------------------
# 1) CORE COMPETENCIES: collapse the three long bullet paragraphs into a
#    single short summary line.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Research and Analytics: Survey Methodology: Design, sampling, weighting, longitudinal analysis $bu Statistical Analysis: Regression modeling, clustering, segmentation, machine learning $bu Geospatial Analysis: Spatial clustering, boundary estimation, demographic mapping $bu Data Visualization: Tableau, PowerBI, d3.js, Matplotlib, Seaborn, choropleth mapping $bu Research Management: Team leadership, methodology design, stakeholder communication",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Research and Analytics $bu Programming and Development $bu Data Infrastructure",
    2) | Out-Null

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Programming and Development: Python:*") {
        $p.Range.Delete()
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Data Infrastructure: Cloud Platforms:*") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 2) RESEARCH DIRECTOR (Progressive Change Campaign Committee) bullets
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Managed critical research operations for political campaigns*") {
        $p.Range.Text = "$bu Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls"
        $p.Range.InsertAfter("$cr$bu Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren")
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Conducted comprehensive polling and demographic analysis*") {
        $p.Range.Text = "$bu Built tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver"
        $p.Range.InsertAfter("$cr$bu Designed survey deployment system facilitating thousands of simultaneous phone surveys")
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Developed strategic recommendations based on data analysis*") {
        $p.Range.Text = "$bu Significantly increased data collection efficiency through automated calling infrastructure"
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Led research team in support of progressive political initiatives*") {
        $p.Range.Text = "$bu Managed comprehensive research operations for progressive political initiatives and candidates"
        break
    }
}

# ---------------------------------------------------------------------
# 3) SOFTWARE ENGINEER (Salsa Labs) bullets
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Developed software solutions for political campaigns and advocacy groups*") {
        $p.Range.Text = "$bu Maintained and extended entire geospatial analysis and reporting tools for Java-based CRM system"
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Built web applications for voter engagement and campaign management*") {
        $p.Range.Text = "$bu Developed custom tile server for Web Map Service (WMS) integration using GeoTools and OpenLayers"
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Integrated third-party APIs and data sources for campaign tools*") {
        $p.Range.Text = "$bu Built geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill"
        $p.Range.InsertAfter("$cr$bu Integrated mapping and visualization tools for political campaign data analysis")
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Collaborated with political strategists to translate requirements into technical solutions*") {
        $p.Range.Text = "$bu Collaborated with political strategists to translate geospatial requirements into technical solutions"
        break
    }
}

# ---------------------------------------------------------------------
# 4) INTERIM TECHNOLOGY MANAGER (The Praxis Project) bullets
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Integrated technology solutions within organizational frameworks for social justice organizations*") {
        $p.Range.Text = "$bu Assisted in search for full-time CTO while performing all programmatic technology roles for multi-million dollar organization"
        $p.Range.InsertAfter("$cr$bu Made all technology decisions and practices for massive multinational non-governmental organization")
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Developed data management systems for community organizing efforts*") {
        $p.Range.Text = "$bu Wrote comprehensive frameworks for internal and external technology audits"
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Provided technical training and support to nonprofit staff*") {
        $p.Range.Text = "$bu Trained beneficiaries on spatial and Census data analysis for public health research"
        $p.Range.InsertAfter("$cr$bu Trained NGO staff in web development using Drupal, PHP, and MySQL")
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Built custom applications for community engagement and advocacy*") {
        $p.Range.Text = "$bu Managed technology infrastructure supporting community health initiatives across multiple countries"
        break
    }
}

# ---------------------------------------------------------------------
# 5) PROGRAMMER (Lake Research Partners) bullets
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Developed data analysis tools for political polling and research*") {
        $p.Range.Text = "$bu Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party"
        $p.Range.InsertAfter("$cr$bu Developed system that later became the Polling Consortium Database at The Analyst Institute")
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Built statistical models for voter behavior analysis*") {
        $p.Range.Text = "$bu Worked on all aspects of questionnaire design, sampling, reporting and analysis for Congressional, Senate and Presidential elections"
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Created data visualization tools for research presentations*") {
        $p.Range.Text = "$bu Conducted statistical modeling and analysis using SPSS, ArcGIS, Quantum GIS, GRASS, Stata, OSCAR, PostgreSQL, PostGIS, and Oracle"
        $p.Range.InsertAfter("$cr$bu Pioneered integration of advanced mapping techniques into standard reports including choropleths and hexagonal grid maps")
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Supported senior researchers with technical analysis and reporting*") {
        $p.Range.Text = "$bu Developed innovative approaches to visualizing demographic and market data for enhanced client understanding"
        break
    }
}

# ---------------------------------------------------------------------
# 6) FIELD DIRECTOR (The Feldman Group) bullets
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Managed field operations for political campaigns and research projects*") {
        $p.Range.Text = "$bu Administered all quantitative and qualitative research operations ensuring reporting accuracy"
        $p.Range.InsertAfter("$cr$bu Managed comprehensive survey fielding for multi-million dollar research firm")
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Developed data collection and management systems for field work*") {
        $p.Range.Text = "$bu Developed and implemented data warehousing solutions for efficient storage and retrieval of research findings"
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Trained field staff on data collection protocols and quality control*") {
        $p.Range.Text = "$bu Created custom reports and data visualizations based on specific client requirements"
        $p.Range.InsertAfter("$cr$bu Introduced mapping and geospatial analysis into standard reporting procedures")
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Analyzed field data to inform campaign strategy and research findings*") {
        $p.Range.Text = "$bu Enhanced value of research deliverables through advanced analytical techniques using SPSS, OSCAR, PHP, and MySQL"
        break
    }
}

# ---------------------------------------------------------------------
# 7) New TECHNICAL SKILLS section at the end of the document.
#    Insert all paragraphs first (while still "Normal" style), then
#    style the heading last so the body lines below it are not affected
#    by style inheritance from the heading paragraph.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertAfter("${cr}TECHNICAL SKILLS")
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertAfter("${cr}RESEARCH AND ANALYTICS Survey Methodology; Statistical Analysis; Geospatial Analysis; Data Visualization; Research Management")
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertAfter("${cr}PROGRAMMING AND DEVELOPMENT Python; JVM Languages; Web Technologies; Database Languages; Statistical Computing")
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertAfter("${cr}DATA INFRASTRUCTURE Cloud Platforms; Big Data; Databases; Geospatial; DevOps")

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*TECHNICAL SKILLS*") {
        $p.Style = "Heading2"
        break
    }
}

Write-Host "Done"
